$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells K1/L1 ("volumen" / "flujo")
$ws.Range("K1").Value = "volumen"
$ws.Range("L1").Value = "flujo"

# Fill data rows 2-20 for new columns K (volumen=200) and L (flujo=12)
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 11).Value = 200
    $ws.Cells.Item($r, 12).Value = 12
}

# Match the number-format style used by column C (s="1", numFmtId 49 "text")
$ws.Range("K2:L20").NumberFormat = "@"

# Update selection to match the diff (K5 active cell)
$ws.Range("K5").Select()
